$d = $word.ActiveDocument

# --- 1. Split the "bcrypt" paragraph's single run into three runs, wrapping
#        "sécurité" with spell-check proofErr markers (as Word's
#        type-as-you-go spell checker would do). ---
$bcryptIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "bcrypt :*") {
        $bcryptIndex = $i
    }
}

$bcryptPara = $d.Paragraphs($bcryptIndex)
$bcryptXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">bcrypt : Permet de hasher les mots de passe pour plus de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>s&#233;curit&#233;</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>'
$bcryptPara.Range.InsertXML($bcryptXml)

# --- 2. Add a brand-new paragraph after it: "npm install cookie-parser",
#        with "npm" wrapped in spell-check proofErr markers. ---
$bcryptPara2 = $d.Paragraphs($bcryptIndex)
$bcryptPara2.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($bcryptIndex + 1)
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> install cookie-parser</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newXml)
